# GroupDocs.Assembly 25.12 December Christmas Release update
#
# 1. Bump the evaluation-watermark version text that is visible in the
#    document body from "25.6" to "25.12".
# 2. Register the standard Word "Hyperlink" character style in the
#    styles part (it is referenced/expected by newer GroupDocs.Assembly
#    output but was missing from this document's style table).

$d = $word.ActiveDocument

# --- 1. Update the visible evaluation watermark text -----------------
$d.Content.Find.Execute(
    "Evaluation Only. Created with GroupDocs.Assembly 25.6. " + [char]0x00A9 + " Aspose Pty Ltd 2001-2025. All Rights Reserved.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Evaluation Only. Created with GroupDocs.Assembly 25.12. " + [char]0x00A9 + " Aspose Pty Ltd 2001-2025. All Rights Reserved.",
    2
) | Out-Null

# --- 2. Add the built-in "Hyperlink" character style ------------------
$hyperlinkStyle = $d.Styles.Add("Hyperlink", 2)
$hyperlinkStyle.BaseStyle = "DefaultParagraphFont"
$hyperlinkStyle.Priority = 99
$hyperlinkStyle.Font.Color = 12673797
$hyperlinkStyle.Font.Underline = 1
